$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename column C header and add new column D header ---
$ws.Range("C1").Value = "need_to_be_run_with_LEXTREME"
$ws.Range("D1").Value = "need_to_be_run_with_LexGlue"

# --- Data rows 2-35: fill boolean flags for LEXTREME (C) / LexGlue (D) ---
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = $false
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = $false
$ws.Range("C5").Value = $true
$ws.Range("D5").Value = $false
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = $false
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = $true
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $true
$ws.Range("C9").Value = $true
$ws.Range("D9").Value = $false
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = $false
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = $false
$ws.Range("C12").Value = $false
$ws.Range("D12").Value = $false
$ws.Range("C13").Value = $true
$ws.Range("D13").Value = $false
$ws.Range("C14").Value = $false
$ws.Range("D14").Value = $false
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = $false
$ws.Range("C16").Value = $true
$ws.Range("D16").Value = $false
$ws.Range("C17").Value = $true
$ws.Range("D17").Value = $false
$ws.Range("C18").Value = $true
$ws.Range("D18").Value = $false
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = $false
$ws.Range("C20").Value = $true
$ws.Range("D20").Value = $false
$ws.Range("C21").Value = $true
$ws.Range("D21").Value = $false
$ws.Range("C22").Value = $true
$ws.Range("D22").Value = $false
$ws.Range("C23").Value = $true
$ws.Range("D23").Value = $false
$ws.Range("C24").Value = $true
$ws.Range("D24").Value = $false
$ws.Range("C25").Value = $false
$ws.Range("D25").Value = $false
$ws.Range("C26").Value = $true
$ws.Range("D26").Value = $false
$ws.Range("C27").Value = $true
$ws.Range("D27").Value = $false
$ws.Range("C28").Value = $true
$ws.Range("D28").Value = $false
$ws.Range("C29").Value = $true
$ws.Range("D29").Value = $false
$ws.Range("C30").Value = $false
$ws.Range("D30").Value = $false
$ws.Range("C31").Value = $true
$ws.Range("D31").Value = $false
$ws.Range("C32").Value = $true
$ws.Range("D32").Value = $false
$ws.Range("C33").Value = $false
$ws.Range("D33").Value = $false
$ws.Range("C34").Value = $true
$ws.Range("D34").Value = $true
$ws.Range("C35").Value = $true
$ws.Range("D35").Value = $true

# --- Remove the AutoFilter (dropdown arrows) that covered A1:C35 ---
$ws.AutoFilterMode = $false

# --- Update the hidden _FilterDatabase defined name to the new extent A1:D35 ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Tabelle1!_FilterDatabase") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$D`$35"
    }
}

# --- Update the view: zoom + selected cell to match author session ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 136
$ws.Range("D37").Select()

Write-Output "done"
